$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (single decimal point)
# must be forced to Text format first, so Excel keeps them as text instead of
# auto-converting them to numbers (matching the source data which stores these as strings).
$textCells = @("D5", "D6", "D11", "D16", "D20", "D25", "D27", "D34", "D36", "D45", "D48", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.986.75'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.555.15'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = '207.18'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '0.486'
$ws.Range("E6").Value = '  +1.09%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '1.777.66'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '1.556.51'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '61.93'
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.978.21'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '7.28'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").Value = '152.44'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("D27").Value = '14.91'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").Value = '1.402.39'
$ws.Range("E33").Value = '  +4.85%  '
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +2.97%  '
$ws.Range("E35").Value = '  +3.39%  '
$ws.Range("D36").Value = '0.953'
$ws.Range("E36").Value = '  +2.76%  '
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("E44").Value = '  -4.38%  '
$ws.Range("D45").Value = '63.86'
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").Value = '1.691.06'
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").Value = '86.25'
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("D50").Value = '0.0956'
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("E51").Value = '  +0.64%  '
